$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by cloning the existing "2022-Q3"
#    sheet (item 2) so that formatting/styles (header row, column A
#    style, etc.) come along for free, then placing the clone right
#    before it.
# ------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item(2)
$sheetQ3.Copy($sheetQ3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Trim the copied sheet from 7 data rows down to the 5 needed for
# 2022-Q4 (rows 7..9 held the extra funds from the Q3 sheet).
$newSheet.Rows.Item(7).Resize(3).Delete()

# Force columns B and D:G to be stored as text (the source data keeps
# these as formatted strings - e.g. fund codes like "004702" would
# otherwise lose their leading zero, and "12.97" would become a float).
$newSheet.Range("B2:B6").NumberFormat = "@"
$newSheet.Range("D2:G6").NumberFormat = "@"


# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "004702"
$newSheet.Range("C2").Value = "南方金融主题灵活配置混合A"
$newSheet.Range("D2").Value = "12.97"
$newSheet.Range("E2").Value = "92.71"
$newSheet.Range("F2").Value = "3.57"
$newSheet.Range("G2").Value = "0.4630"
$newSheet.Range("H2").Value = 9

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "013500"
$newSheet.Range("C3").Value = "南方金融主题灵活配置混合C"
$newSheet.Range("D3").Value = "4.80"
$newSheet.Range("E3").Value = "92.71"
$newSheet.Range("F3").Value = "3.57"
$newSheet.Range("G3").Value = "0.1714"
$newSheet.Range("H3").Value = 9

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "001244"
$newSheet.Range("C4").Value = "华泰柏瑞量化智慧灵活配置混合A"
$newSheet.Range("D4").Value = "2.89"
$newSheet.Range("E4").Value = "93.57"
$newSheet.Range("F4").Value = "0.67"
$newSheet.Range("G4").Value = "0.0194"
$newSheet.Range("H4").Value = 3

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "516980"
$newSheet.Range("C5").Value = "华富中证证券公司先锋策略ETF"
$newSheet.Range("D5").Value = "0.25"
$newSheet.Range("E5").Value = "99.41"
$newSheet.Range("F5").Value = "4.26"
$newSheet.Range("G5").Value = "0.0106"
$newSheet.Range("H5").Value = 6

# Row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "006104"
$newSheet.Range("C6").Value = "华泰柏瑞量化智慧灵活配置混合C"
$newSheet.Range("D6").Value = "0.38"
$newSheet.Range("E6").Value = "93.57"
$newSheet.Range("F6").Value = "0.67"
$newSheet.Range("G6").Value = "0.0025"
$newSheet.Range("H6").Value = 3

# The NumberFormat="@" assignment above stamped a new cell style onto
# B2:B6/D2:G6 (needed so the numeric-looking strings stick as text
# instead of being parsed into numbers). Paste just the *formatting*
# of a never-touched, default-styled cell back over them so the
# visual/style result matches the rest of the (unstyled) data cells.
$newSheet.Range("A1").Copy()
$newSheet.Range("B2:B6").PasteSpecial(-4122)
$newSheet.Range("D2:G6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 for
#    2022-Q4 and push the rest of the quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Capture existing B:D values (quarter label, count, market value)
# for rows 2..8 before they get overwritten, shifting bottom-up.
$rowsData = @()
for ($r = 2; $r -le 8; $r++) {
    $rowsData += , @($summary.Range("B$r").Value2, $summary.Range("C$r").Value2, $summary.Range("D$r").Value2)
}

# Give row 9 the same look (border/bold A-column style) as row 8.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write old rows 2..8 into rows 3..9 (bottom-up so we never clobber
# a value before it has been read... values were already captured
# above, so plain top-down writes are fine too).
for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $destRow = $i + 3
    $summary.Range("A$destRow").Value = $i + 1
    $summary.Range("B$destRow").Value = $rowsData[$i][0]
    $summary.Range("C$destRow").Value = $rowsData[$i][1]
    $summary.Range("D$destRow").Value = $rowsData[$i][2]
}

# New row 2: 2022-Q4 data
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.67

# ------------------------------------------------------------------
# 3. Restore the active sheet/tab to the summary sheet so the
#    workbook view doesn't change unexpectedly.
# ------------------------------------------------------------------
$summary.Activate()

Write-Host "Added 2022-Q4 sheet and updated 总计 summary."
